$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 783, shifting existing rows 783:824 down to 784:825
$ws.Rows.Item(783).Insert()

# Write the date as literal text (not auto-converted to a date serial number):
# force a text number format before assigning, then clear the format again so
# the cell ends up with no explicit style, matching the surrounding data rows.
$dateCell = $ws.Cells.Item(783, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/05"
$dateCell.ClearFormats()

# Remaining columns for the newly inserted row
$ws.Cells.Item(783, 2).Value = "木"
$ws.Cells.Item(783, 3).Value = 5
$ws.Cells.Item(783, 4).Value = 201
